$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cells C2:F6 previously stored numbers as text because the decimal
# point had been replaced by a comma while writing the file (see commit
# message). This re-enters the values as proper numbers; the one value
# that could not be recovered (bbp / 2016Q3) is written as "NA".

$ws.Range("C2").Value = 6.83
$ws.Range("D2").Value = 6.8
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = 6.71

$ws.Range("C3").Value = 23.34
$ws.Range("D3").Value = 23.34
$ws.Range("E3").Value = 23.2
$ws.Range("F3").Value = 23.68

$ws.Range("C4").Value = 19.35
$ws.Range("D4").Value = 19.38
$ws.Range("E4").Value = 19.58
$ws.Range("F4").Value = 19.68

$ws.Range("C5").Value = 3.84
$ws.Range("D5").Value = 3.84
$ws.Range("E5").Value = 3.54
$ws.Range("F5").Value = 3.74

$ws.Range("C6").Value = 19.38
$ws.Range("D6").Value = 19.58
$ws.Range("E6").Value = 19.35
$ws.Range("F6").Value = 19.69

# Give the re-entered numeric block its own explicit black font (distinct
# from the header's style) while keeping the right alignment it already had.
$numRange = $ws.Range("C2:F6")
$numRange.Font.Name = "Calibri"
$numRange.Font.Color = 0

# Reflect the new selection left behind after re-entering the data.
$ws.Range("C2:F6").Select() | Out-Null
